$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("formulário de detalhamento")
$ws.Range("H100").Value = "Test value"
